$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.08661466666666667
$ws.Range("H2").Value = 0.259844
$ws.Range("I2").Value = 0.5374496355558498
$ws.Range("J2").Value = 0.5374496355558498
$ws.Range("Q2").Value = 0.0005500320048888889
$ws.Range("R2").Value = 0.004950288044
$ws.Range("S2").Value = 0.5374496355558498
$ws.Range("T2").Value = 0.5374496355558498

# Row 3
$ws.Range("I3").Value = 0.2407068810034004
$ws.Range("J3").Value = 0.2407068810034003
$ws.Range("S3").Value = 0.2407068810034004
$ws.Range("T3").Value = 0.2407068810034003

# Row 4
$ws.Range("G4").Value = 0.035752
$ws.Range("H4").Value = 0.107256
$ws.Range("I4").Value = 0.2218434834407499
$ws.Range("J4").Value = 0.2218434834407499
$ws.Range("Q4").Value = 0.0002270371173333333
$ws.Range("R4").Value = 0.002043334056
$ws.Range("S4").Value = 0.2218434834407499
$ws.Range("T4").Value = 0.2218434834407499
